$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" column (E16:E28): data reordered from descending to ascending
$periodos = @{
    16 = "2304"
    17 = "2305"
    18 = "2306"
    19 = "2307"
    20 = "2308"
    21 = "2309"
    22 = "2310"
    23 = "2311"
    24 = "2312"
    25 = "2401"
    26 = "2402"
    27 = "2403"
    28 = "2404"
}

foreach ($row in $periodos.Keys) {
    $ws.Range("E$row").Value = $periodos[$row]
}

# Update "Salario Basico" values that moved along with the reordering
$ws.Range("F16").Value = 27840
$ws.Range("F28").Value = 46400

# Update Salario Basico for the last worker row
$ws.Range("G29").Value = 1300000
